$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 121, shifting existing rows 121..225 down to 122..226
$ws.Rows.Item(121).Insert()

# Populate the new row 121 with the new data record
$ws.Cells.Item(121, 1).Value2 = 4
$ws.Cells.Item(121, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(121, 3).Value2 = "Los Lagos"
$ws.Cells.Item(121, 4).Value2 = 44589
$ws.Cells.Item(121, 5).Value2 = 10
$ws.Cells.Item(121, 6).Value2 = 100112040
$ws.Cells.Item(121, 7).Value2 = "Cilantro"
$ws.Cells.Item(121, 8).Value2 = "Sin especificar"
$ws.Cells.Item(121, 9).Value2 = "Primera"
$ws.Cells.Item(121, 10).Value2 = 220
$ws.Cells.Item(121, 11).Value2 = 16000
$ws.Cells.Item(121, 12).Value2 = 16000
$ws.Cells.Item(121, 13).Value2 = 16000
$ws.Cells.Item(121, 14).Value2 = "$/caja 36 atados"
$ws.Cells.Item(121, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(121, 16).Value2 = 444
$ws.Cells.Item(121, 17).Value2 = 36
$ws.Cells.Item(121, 18).Value2 = "Hortaliza"
